# Applies the update described by the commit:
#  - Update "WATER FOR INJECTION AMP. 5 ML" row: balance 8521:0 -> 8520:0,
#    selling price 6.0000 -> 8.0000, transactions 3:0 -> 4:0
#  - Update "سرنجات 3 سم" row: selling price 4.0000 -> 6.0000, transactions 2:0 -> 3:0
#  - Insert a new item row "سرنجات 5 سم" right after it (#41), pushing the
#    remaining rows (كالونا, كريم فيبكس الازرق, totals, footer) down by one
#  - Update the grand total and re-number the shifted rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) WATER FOR INJECTION AMP. 5 ML  (row 42)
# ---------------------------------------------------------------------
$ws.Range("H42").Value = "8520:0"
$ws.Range("P42").Value = "8.0000"
$ws.Range("Q42").Value = "4:0"

# ---------------------------------------------------------------------
# 2) سرنجات 3 سم  (row 46)
# ---------------------------------------------------------------------
$ws.Range("P46").Value = "6.0000"
$ws.Range("Q46").Value = "3:0"

# ---------------------------------------------------------------------
# 3) Insert a new row for "سرنجات 5 سم" after row 46, before the old
#    row 47 (كالونا). Copy row 46's formatting (same item-row style) so
#    the new row matches the look of the surrounding rows, then fill in
#    the new values.
# ---------------------------------------------------------------------
$ws.Rows.Item(47).Insert()
$ws.Range("A46:Q46").Copy()
$ws.Range("A47:Q47").PasteSpecial(-4104)

# restore the thin bottom border lost by the copy/paste (matches the
# other item rows' border)
$newRowBorder = $ws.Range("A47:Q47")
$newRowBorder.Borders.Item(9).LineStyle = 1
$newRowBorder.Borders.Item(9).Color = 13882323

$ws.Range("A47").Value = 41
$ws.Range("C47").Value = "سرنجات 5 سم"
$ws.Range("H47").Value = "0:0"
$ws.Range("L47").Value = "0"
$ws.Range("N47").Value = "3.00"
$ws.Range("P47").Value = "3.0000"
$ws.Range("Q47").Value = "1:0"

# ---------------------------------------------------------------------
# 4) Update the grand total (now shifted one row down to row 50)
# ---------------------------------------------------------------------
$ws.Range("P50").Value = 1842.7750000000001
